$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D2:D51) stays formatted as Text so that values
# like "239.54" or "1.00" are not auto-converted into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "43.936.19"
$ws.Range("E2").Value = "  -0.59%  "

# Row 3
$ws.Range("D3").Value = "2.347.67"
$ws.Range("E3").Value = "  -1.22%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "239.54"
$ws.Range("E5").Value = "  -1.66%  "

# Row 6
$ws.Range("D6").Value = "0.669"
$ws.Range("E6").Value = "  -3.62%  "

# Row 7
$ws.Range("D7").Value = "71.99"
$ws.Range("E7").Value = "  -7.21%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -2.57%  "

# Row 10
$ws.Range("D10").Value = "0.0997"
$ws.Range("E10").Value = "  -4.76%  "

# Row 11
$ws.Range("D11").Value = "58.10"
$ws.Range("E11").Value = "  +0.40%  "

# Row 12
$ws.Range("D12").Value = "32.54"
$ws.Range("E12").Value = "  -0.20%  "

# Row 13
$ws.Range("E13").Value = "  -0.48%  "

# Row 14
$ws.Range("D14").Value = "7.21"
$ws.Range("E14").Value = "  -4.63%  "

# Row 15
$ws.Range("D15").Value = "2.699.25"
$ws.Range("E15").Value = "  -1.25%  "

# Row 16
$ws.Range("D16").Value = "16.27"
$ws.Range("E16").Value = "  -5.71%  "

# Row 17
$ws.Range("D17").Value = "0.900"
$ws.Range("E17").Value = "  -3.24%  "

# Row 18
$ws.Range("D18").Value = "2.352.25"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19
$ws.Range("D19").Value = "43.812.26"
$ws.Range("E19").Value = "  -1.39%  "

# Row 20
$ws.Range("D20").Value = "0.0000102"
$ws.Range("E20").Value = "  -1.86%  "

# Row 21
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").Value = "  -0.79%  "

# Row 22
$ws.Range("D22").Value = "78.19"
$ws.Range("E22").Value = "  -0.95%  "

# Row 23
$ws.Range("D23").Value = "253.41"
$ws.Range("E23").Value = "  -1.87%  "

# Row 24
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  +6.04%  "

# Row 25
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("E26").Value = "  +0.52%  "

# Row 27
$ws.Range("E27").Value = "  -2.96%  "

# Row 28
$ws.Range("D28").Value = "10.42"
$ws.Range("E28").Value = "  -5.14%  "

# Row 29
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  +0.44%  "

# Row 30
$ws.Range("D30").Value = "177.32"
$ws.Range("E30").Value = "  +0.85%  "

# Row 31
$ws.Range("D31").Value = "22.27"
$ws.Range("E31").Value = "  -3.81%  "

# Row 32
$ws.Range("E32").Value = "  -2.73%  "

# Row 33
$ws.Range("D33").Value = "0.136"
$ws.Range("E33").Value = "  +0.16%  "

# Row 34
$ws.Range("D34").Value = "0.0745"
$ws.Range("E34").Value = "  -2.67%  "

# Row 35
$ws.Range("D35").Value = "5.11"
$ws.Range("E35").Value = "  -5.24%  "

# Row 36
$ws.Range("D36").Value = "5.38"
$ws.Range("E36").Value = "  +0.42%  "

# Row 37
$ws.Range("D37").Value = "3.74"
$ws.Range("E37").Value = "  -4.29%  "

# Row 38
$ws.Range("D38").Value = "6.41"
$ws.Range("E38").Value = "  -3.05%  "

# Row 39
$ws.Range("D39").Value = "2.37"
$ws.Range("E39").Value = "  -5.28%  "

# Row 40
$ws.Range("D40").Value = "0.0274"
$ws.Range("E40").Value = "  -1.78%  "

# Row 41
$ws.Range("D41").Value = "66.36"
$ws.Range("E41").Value = "  +20.65%  "

# Row 42
$ws.Range("D42").Value = "5.28"
$ws.Range("E42").Value = "  +17.64%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "9.17"
$ws.Range("E43").Value = "  -0.18%  "

# Row 44
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.108"
$ws.Range("E44").Value = "  +6.43%  "

# Row 45
$ws.Range("D45").Value = "18.71"
$ws.Range("E45").Value = "  -2.44%  "

# Row 46
$ws.Range("D46").Value = "0.198"
$ws.Range("E46").Value = "  +0.29%  "

# Row 47
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("E48").Value = "  -3.32%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.46"
$ws.Range("E49").Value = "  -3.14%  "

# Row 50
$ws.Range("D50").Value = "98.74"
$ws.Range("E50").Value = "  -4.37%  "

# Row 51
$ws.Range("E51").Value = "  -6.19%  "
